$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column G ("hire_date") into column H: same header text and the
# same hire-date values/format for every data row.
$ws.Range("H1").Value = $ws.Range("G1").Value2

$lastRow = 16
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 8).NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
}

# New column H is a bit wider than the old column G already was.
$ws.Columns.Item(8).ColumnWidth = 19.6

# Move the active selection to J13, matching the saved selection state.
$ws.Range("J13").Select() | Out-Null
